# Generate Report for Handback
#
# Populates the "Latest HO Xliff Generate Date" / "Correspond Handoff
# Datetime" / "Correspond Handback DateTime" columns for the second
# tracked file (5702a499-1e12-4816-b4f0-8e0f0822bf6e) now that its
# handback report has been generated. The first file
# (409ba8ab-6a1c-4494-9e24-d3cd1508d24c) already had these values and is
# left untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# zh-cn handoff/handback timestamps for row 3 (5702a499...)
$zhcn.Range("H3").Value = "2016-09-04 10:53:08"
$zhcn.Range("K3").Value = "2016-09-04 10:53:34"

# de-de handoff/handback timestamps for row 3 (5702a499...)
$dede.Range("H3").Value = "2016-09-04 10:53:13"
$dede.Range("K3").Value = "2016-09-04 10:53:41"

# Overview sheet: latest HO xliff generate date for that file is the
# newest of its per-language handoff dates (de-de's 10:53:13).
$overview.Range("G3").Value = "2016-09-04 10:53:13"
